$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 220111.8
$ws.Cells.Item(6, 9).Value = 250140
$ws.Cells.Item(6, 11).Value = 750420
$ws.Cells.Item(6, 13).Value = -750308

$ws.Cells.Item(9, 8).Value = 207.41176
$ws.Cells.Item(9, 9).Value = 98.61539
$ws.Cells.Item(9, 11).Value = 98.61539
$ws.Cells.Item(9, 13).Value = 70.38461

$ws.Cells.Item(26, 8).Value = 0
$ws.Cells.Item(26, 10).Value = 0
$ws.Cells.Item(26, 12).Value = 0
$ws.Cells.Item(26, 14).ClearContents()

$ws.Cells.Item(41, 8).Value = 743.0769
$ws.Cells.Item(41, 9).Value = 274.33334
$ws.Cells.Item(41, 11).Value = 274.33334
$ws.Cells.Item(41, 13).Value = 165.66666

$ws.Cells.Item(107, 8).Value = 11230.926
$ws.Cells.Item(107, 9).Value = 10402.087
$ws.Cells.Item(107, 11).Value = 10402.087
$ws.Cells.Item(107, 13).Value = -8482.087

$ws.Cells.Item(113, 8).Value = 4999.4
$ws.Cells.Item(113, 9).Value = 4999.4
$ws.Cells.Item(113, 11).Value = 4999.4
$ws.Cells.Item(113, 13).Value = -1745.4

$ws.Cells.Item(116, 8).Value = 18037.842
$ws.Cells.Item(116, 9).Value = 4072.1
$ws.Cells.Item(116, 10).Value = 33555.332
$ws.Cells.Item(116, 11).Value = 4072.1
$ws.Cells.Item(116, 12).Value = 33555.332
$ws.Cells.Item(116, 13).Value = -630.0999999999999
$ws.Cells.Item(116, 14).Value = -40439.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(37, 8).Value = 20000
$ws.Cells.Item(37, 10).Value = 20000
$ws.Cells.Item(37, 12).Value = 20000
$ws.Cells.Item(37, 14).Value = -20546

$ws.Cells.Item(88, 8).Value = 2853.6667
$ws.Cells.Item(88, 10).Value = 2767.1667
$ws.Cells.Item(88, 12).Value = 2767.1667
$ws.Cells.Item(88, 14).Value = -3579.1667

$ws.Cells.Item(91, 8).Value = 2853.6667
$ws.Cells.Item(91, 10).Value = 2767.1667
$ws.Cells.Item(91, 12).Value = 2767.1667
$ws.Cells.Item(91, 14).Value = -5575.1667

$ws.Cells.Item(102, 8).Value = 3391.3333
$ws.Cells.Item(102, 9).Value = 3789
$ws.Cells.Item(102, 11).Value = 3789
$ws.Cells.Item(102, 13).Value = -2167

$ws.Cells.Item(122, 8).Value = 21740682
$ws.Cells.Item(122, 9).Value = 31251316
$ws.Cells.Item(122, 10).Value = 2091.8572
$ws.Cells.Item(122, 11).Value = 93753948
$ws.Cells.Item(122, 12).Value = 6275.571599999999
$ws.Cells.Item(122, 13).Value = -93751498
$ws.Cells.Item(122, 14).Value = -11175.5716

$ws.Cells.Item(135, 8).Value = 79999
$ws.Cells.Item(135, 10).Value = 79999
$ws.Cells.Item(135, 12).Value = 79999
$ws.Cells.Item(135, 14).Value = -90139

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(63, 8).Value = 60000
$ws.Cells.Item(63, 10).Value = 60000
$ws.Cells.Item(63, 12).Value = 60000
$ws.Cells.Item(63, 14).Value = -61372

$ws.Cells.Item(66, 8).Value = 60000
$ws.Cells.Item(66, 10).Value = 60000
$ws.Cells.Item(66, 12).Value = 180000
$ws.Cells.Item(66, 14).Value = -186864

$ws.Cells.Item(99, 8).Value = 6767.107
$ws.Cells.Item(99, 9).Value = 7803.4346
$ws.Cells.Item(99, 11).Value = 7803.4346
$ws.Cells.Item(99, 13).Value = -6305.4346

$ws.Cells.Item(134, 8).Value = 3211573.8
$ws.Cells.Item(134, 9).Value = 4680.6045
$ws.Cells.Item(134, 10).Value = 18533396
$ws.Cells.Item(134, 11).Value = 14041.8135
$ws.Cells.Item(134, 12).Value = 55600188
$ws.Cells.Item(134, 13).Value = -11506.8135
$ws.Cells.Item(134, 14).Value = -55605258

$ws.Cells.Item(135, 8).Value = 99999
$ws.Cells.Item(135, 10).Value = 99999
$ws.Cells.Item(135, 12).Value = 99999
$ws.Cells.Item(135, 14).Value = -110139

$ws.Cells.Item(137, 8).Value = 78332.336
$ws.Cells.Item(137, 10).Value = 78332.336
$ws.Cells.Item(137, 12).Value = 78332.336
$ws.Cells.Item(137, 14).Value = -88532.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 62577
$ws.Cells.Item(7, 9).Value = 90946
$ws.Cells.Item(7, 10).Value = 165.2
$ws.Cells.Item(7, 11).Value = 90946
$ws.Cells.Item(7, 12).Value = 165.2
$ws.Cells.Item(7, 13).Value = -90833
$ws.Cells.Item(7, 14).Value = -391.2

$ws.Cells.Item(10, 8).Value = 2794
$ws.Cells.Item(10, 9).Value = 1352.8
$ws.Cells.Item(10, 11).Value = 1352.8
$ws.Cells.Item(10, 13).Value = -1213.8

$ws.Cells.Item(31, 8).Value = 1324882.6
$ws.Cells.Item(31, 9).Value = 1738197.4
$ws.Cells.Item(31, 10).Value = 2275.2
$ws.Cells.Item(31, 11).Value = 1738197.4
$ws.Cells.Item(31, 12).Value = 2275.2
$ws.Cells.Item(31, 13).Value = -1737902.4
$ws.Cells.Item(31, 14).Value = -2865.2

$ws.Cells.Item(34, 8).Value = 1324882.6
$ws.Cells.Item(34, 9).Value = 1738197.4
$ws.Cells.Item(34, 10).Value = 2275.2
$ws.Cells.Item(34, 11).Value = 1738197.4
$ws.Cells.Item(34, 12).Value = 2275.2
$ws.Cells.Item(34, 13).Value = -1737995.4
$ws.Cells.Item(34, 14).Value = -2679.2

$ws.Cells.Item(86, 8).Value = 11215.4
$ws.Cells.Item(86, 9).Value = 17927.285
$ws.Cells.Item(86, 11).Value = 17927.285
$ws.Cells.Item(86, 13).Value = -16804.285

$ws.Cells.Item(89, 8).Value = 11215.4
$ws.Cells.Item(89, 9).Value = 17927.285
$ws.Cells.Item(89, 11).Value = 89636.425
$ws.Cells.Item(89, 13).Value = -84020.425

$ws.Cells.Item(122, 8).Value = 9767.190000000001
$ws.Cells.Item(122, 9).Value = 1829.8125
$ws.Cells.Item(122, 10).Value = 35166.8
$ws.Cells.Item(122, 11).Value = 5489.4375
$ws.Cells.Item(122, 12).Value = 105500.4
$ws.Cells.Item(122, 13).Value = -3039.4375
$ws.Cells.Item(122, 14).Value = -110400.4

$ws.Cells.Item(132, 8).Value = 2834.84
$ws.Cells.Item(132, 9).Value = 2963.0588
$ws.Cells.Item(132, 11).Value = 8889.1764
$ws.Cells.Item(132, 13).Value = -6359.1764

$ws.Cells.Item(134, 8).Value = 1729.9615
$ws.Cells.Item(134, 9).Value = 1639.16
$ws.Cells.Item(134, 11).Value = 4917.48
$ws.Cells.Item(134, 13).Value = -2382.48

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(15, 8).Value = 520
$ws.Cells.Item(15, 9).Value = 750
$ws.Cells.Item(15, 11).Value = 2250
$ws.Cells.Item(15, 13).Value = -2110

$ws.Cells.Item(107, 8).Value = 351.1111
$ws.Cells.Item(107, 9).Value = 435.5
$ws.Cells.Item(107, 11).Value = 1306.5
$ws.Cells.Item(107, 13).Value = 613.5

$ws.Cells.Item(117, 8).Value = 6896.636
$ws.Cells.Item(117, 9).Value = 1208
$ws.Cells.Item(117, 10).Value = 10147.286
$ws.Cells.Item(117, 11).Value = 3624
$ws.Cells.Item(117, 12).Value = 30441.858
$ws.Cells.Item(117, 13).Value = -182
$ws.Cells.Item(117, 14).Value = -37325.858

$ws.Cells.Item(131, 8).Value = 1955.8235
$ws.Cells.Item(131, 9).Value = 1451
$ws.Cells.Item(131, 10).Value = 2881.3333
$ws.Cells.Item(131, 11).Value = 4353
$ws.Cells.Item(131, 12).Value = 8643.999899999999
$ws.Cells.Item(131, 13).Value = 687
$ws.Cells.Item(131, 14).Value = -18723.9999

$ws.Cells.Item(134, 8).Value = 5408
$ws.Cells.Item(134, 9).Value = 2745.7778
$ws.Cells.Item(134, 11).Value = 8237.3334
$ws.Cells.Item(134, 13).Value = -3167.3334

$ws.Cells.Item(136, 8).Value = 7241.091
$ws.Cells.Item(136, 9).Value = 5517.1113
$ws.Cells.Item(136, 11).Value = 16551.3339
$ws.Cells.Item(136, 13).Value = -11451.3339

$ws.Cells.Item(138, 8).Value = 36972.92
$ws.Cells.Item(138, 9).Value = 36972.92
$ws.Cells.Item(138, 11).Value = 110918.76
$ws.Cells.Item(138, 13).Value = -105778.76

$ws.Cells.Item(139, 8).Value = 4371.381
$ws.Cells.Item(139, 9).Value = 3160
$ws.Cells.Item(139, 10).Value = 7399.8335
$ws.Cells.Item(139, 11).Value = 9480
$ws.Cells.Item(139, 12).Value = 22199.5005
$ws.Cells.Item(139, 13).Value = -4340
$ws.Cells.Item(139, 14).Value = -32479.5005

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 4797.478
$ws.Cells.Item(122, 10).Value = 3835.5557
$ws.Cells.Item(122, 12).Value = 11506.6671
$ws.Cells.Item(122, 14).Value = -16406.6671

$ws.Cells.Item(132, 8).Value = 11051.731
$ws.Cells.Item(132, 9).Value = 9331.527
$ws.Cells.Item(132, 10).Value = 23437.2
$ws.Cells.Item(132, 11).Value = 27994.581
$ws.Cells.Item(132, 12).Value = 70311.60000000001
$ws.Cells.Item(132, 13).Value = -25464.581
$ws.Cells.Item(132, 14).Value = -75371.60000000001

$ws.Cells.Item(139, 8).Value = 50000
$ws.Cells.Item(139, 10).Value = 50000
$ws.Cells.Item(139, 12).Value = 50000
$ws.Cells.Item(139, 14).Value = -60280

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5879.8
$ws.Cells.Item(7, 9).Value = 6488.778
$ws.Cells.Item(7, 11).Value = 6488.778
$ws.Cells.Item(7, 13).Value = -6376.778

$ws.Cells.Item(41, 8).Value = 8000
$ws.Cells.Item(41, 9).Value = 8000
$ws.Cells.Item(41, 10).Value = 0
$ws.Cells.Item(41, 11).Value = 8000
$ws.Cells.Item(41, 12).Value = 0
$ws.Cells.Item(41, 13).Value = -7562
$ws.Cells.Item(41, 14).ClearContents()

$ws.Cells.Item(46, 8).Value = 4637.3335
$ws.Cells.Item(46, 9).Value = 1666.6666
$ws.Cells.Item(46, 11).Value = 1666.6666
$ws.Cells.Item(46, 13).Value = -1478.6666

$ws.Cells.Item(61, 8).Value = 15645.75
$ws.Cells.Item(61, 9).Value = 13110.308
$ws.Cells.Item(61, 11).Value = 13110.308
$ws.Cells.Item(61, 13).Value = -12908.308

$ws.Cells.Item(113, 8).Value = 15645.75
$ws.Cells.Item(113, 9).Value = 13110.308
$ws.Cells.Item(113, 11).Value = 13110.308
$ws.Cells.Item(113, 13).Value = -10940.308

$ws.Cells.Item(123, 8).Value = 62000
$ws.Cells.Item(123, 10).Value = 62000
$ws.Cells.Item(123, 12).Value = 62000
$ws.Cells.Item(123, 14).Value = -71800

$ws.Cells.Item(126, 8).Value = 5879.8
$ws.Cells.Item(126, 9).Value = 6488.778
$ws.Cells.Item(126, 11).Value = 19466.334
$ws.Cells.Item(126, 13).Value = -16996.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 2935.4167
$ws.Cells.Item(81, 9).Value = 1338.3334
$ws.Cells.Item(81, 10).Value = 7726.6665
$ws.Cells.Item(81, 11).Value = 2676.6668
$ws.Cells.Item(81, 12).Value = 15453.333
$ws.Cells.Item(81, 13).Value = -1615.6668
$ws.Cells.Item(81, 14).Value = -17575.333

$ws.Cells.Item(84, 8).Value = 2935.4167
$ws.Cells.Item(84, 9).Value = 1338.3334
$ws.Cells.Item(84, 10).Value = 7726.6665
$ws.Cells.Item(84, 11).Value = 13383.334
$ws.Cells.Item(84, 12).Value = 77266.66500000001
$ws.Cells.Item(84, 13).Value = -8079.333999999999
$ws.Cells.Item(84, 14).Value = -87874.66500000001

$ws.Cells.Item(122, 8).Value = 60092.367
$ws.Cells.Item(122, 9).Value = 941.05554
$ws.Cells.Item(122, 11).Value = 2823.16662
$ws.Cells.Item(122, 13).Value = -373.16662

$ws.Cells.Item(135, 8).Value = 78999
$ws.Cells.Item(135, 10).Value = 78999
$ws.Cells.Item(135, 12).Value = 78999
$ws.Cells.Item(135, 14).Value = -89139
